$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values (formulas recalc automatically)
$ws.Range("B2").Value = 150.87
$ws.Range("C2").Value = 169.01

$ws.Range("B3").Value = 154.59
$ws.Range("C3").Value = 127.76

$ws.Range("B8").Value = 28

$ws.Range("B17").Value = 51.55

$ws.Range("B21").Value = 14
$ws.Range("B22").Value = 2646

# Update the selection shown in the sheet view
$ws.Range("B18").Select()
